$d = $word.ActiveDocument

# Helper: given a Range that lies inside some paragraph (already edited to
# hold its final text), insert a brand-new paragraph right after it and
# fill that new paragraph with the given text in italics. The index of
# the new (empty) paragraph is re-fetched through the document's
# Paragraphs collection, since a stray/rebuilt Range handle does not
# reliably track a freshly-inserted paragraph boundary.
function Insert-ItalicParagraphAfter($range, [string]$text) {
    $paraIndex = $range.Paragraphs.Item(1).Index
    $range.Collapse(0)
    $range.InsertParagraphAfter() | Out-Null

    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $ip = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $ip.InsertAfter($text)
    $ip.Italic = 1
}

# ---------------------------------------------------------------------
# 1) Ativação date bump
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Objetivos paragraph: replace text + add new italic (English) paragraph
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Apresentar aos alunos as bases teóricas e experimentais dos métodos instrumentais (quantitativos e qualitativos) de uso mais frequente na área química voltada para os parâmetros das legislações ambientais, incluindo o preparo de amostras e a correta avaliação dos resultados analíticos. Ao final da disciplina, o aluno deve ser capaz de escolher e aplicar a metodologia mais adequada à solução dos problemas analíticos voltados ao meio ambiente.") | Out-Null
$rng.Text = "Geral: Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental. Executar procedimentos de análises volumétricas bem como interpretar, avaliar e criticar os resultados obtidos. Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Conhecer os métodos de preparação de amostras a serem analisadas. Compreender as bases teóricas da química analítica quantitativa de interesse ambiental. Executar procedimentos de análises químicas volumétricas (volumetria por neutralização, volumetria por precipitação, volumetria por oxi-redução e volumetria por complexação), análises ambientais de DBO, DQO, OD e turbidez, bem como interpretar, avaliar e criticar os resultados obtidos."

Insert-ItalicParagraphAfter $rng "General: Train professionals at a higher level with the ability to know the sequence of chemical analysis procedures of environmental interest. Perform volumetric analysis as well as interpret, evaluate and criticize the results obtained.Specifics: Train professionals at a higher level with the ability to know the steps of the analytical sequence. Know the methods of preparation to be analysed. Understand the theoretical bases of quantitative analytical chemistry of environmental interest. Procedures of volumetric agglomerated analysis (volumetric by neutralization, volumetric by back, volumetric by  oxidation-reduction  and volumetric by complexation), environmental analyses of BOD, COD, OD and turbidity, as well as to interpret, evaluate and criticize the results."

# ---------------------------------------------------------------------
# 3) Programa resumido paragraph: replace text + add new italic (English) paragraph
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Introdução à Análise Instrumental. Métodos Espectrofotométricos: UV/Visível. Métodos Espectrométricos: Absorção Atômica, Emissão Atômica. Métodos Eletroanalíticos: Potenciometria e Condutimetria. Métodos Cromatográficos: Cromatografia a Gás e Cromatografia Líquida de Alta Eficiência. Métodos Espectroscópicos: Infravermelho e RMN. Outros Métodos: Fluorescência de Raio X, TOC.") | Out-Null
$rng.Text = "Fundamentos da análise titulométrica (preparação de amostras e cálculos). Análises titulométricas por precipitação, neutralização, complexação e oxirredução. Análises ambientais de DBO, DQO, OD e turbidez."

Insert-ItalicParagraphAfter $rng "Fundamentals of titrometric analysis (preparation of calculations and calculations). Titulometric analyses by exclusive, neutralization, complexation and redox. Environmental analyses of BOD, COD, OD and turbidity."

# ---------------------------------------------------------------------
# 4) Programa paragraph (with several w:br-separated items): replace
#    entire content (collapsing the breaks away) + add new italic
#    (English) paragraph. The paragraph is much longer than the other
#    ones, so Find only locks onto its first sentence; widen the range
#    up to the paragraph's real end (read via the Paragraphs collection)
#    before overwriting the text, keeping everything on the *same*
#    Range object throughout (a hop through Paragraph.Range loses the
#    ability to resolve Paragraphs.Item(1).Index correctly downstream).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Introdução à Análise Instrumental. Correlação entre métodos analíticos instrumentais e por via úmida.") | Out-Null
$progIndex = $rng.Paragraphs.Item(1).Index
$rng.End = $d.Paragraphs.Item($progIndex).Range.End
$rng.Text = "- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.- Análises quantitativas em solo, água, alimentos etc.- Análises ambientais: DBO, DQO, OD e turbidez"

Insert-ItalicParagraphAfter $rng "- Fundamentals of titulometric analysis and calculations in titulometric analysis.- Neutralization titrometry: fundamentals, titration indicators, acid base titration curves.- Complexometric titrometry: fundamentals, complexometry with EDTA.- Titration of oxireduction: fundamentals and main indicators- Precipitation titrometry: fundamentals, indicators, argentimetry.- Quantitative analysis of soil, water, food, etc.- Environmental analysis: BOD, COD, OD and turbidity"

# ---------------------------------------------------------------------
# 5) Avaliação: Método / Critério / Norma de recuperação text swaps
# ---------------------------------------------------------------------
$d.Content.Find.Execute("A avaliação da disciplina será feita por meio de avaliações escritas individuais (provas) e avaliações de atividades em grupo (relatórios das aulas práticas, trabalhos escritos e/ou seminários).", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.", 2) | Out-Null

$d.Content.Find.Execute("A média final (MF) será calculada pela média entre duas avaliações teóricas (individuais) e trabalhos experimentais (grupos), este em função das atividades práticas realizadas durante cada bimestre, sendo as avaliações individuais correspondentes a 75 % da composição de MF e as avaliações em grupo correspondentes a 25 % da composição de MF.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.", 2) | Out-Null

$d.Content.Find.Execute("A Nota de Recuperação (NR) será dada pela média aritmética entre a Nota do Semestre (MF) e a Prova da Recuperação (PR), sendo considerado aprovado o aluno que obtiver NR maior ou igual a cinco.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.", 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Bibliografia paragraph: collapse the whole list into a single run
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Willard, H.H., Merrite, L. e Deab, J., INSTRUMENTAÇÃO ANALÍTICA, Fundação Calouste Gulbekian, Lisboa, 1989.") | Out-Null
$biblioIndex = $rng.Paragraphs.Item(1).Index
$rng.Start = $d.Paragraphs.Item($biblioIndex).Range.Start
$rng.End = $d.Paragraphs.Item($biblioIndex).Range.End
$rng.Text = "Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002."
